# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the f2bb82f3-* row on both the zh-cn and de-de report
# sheets, to reflect newly (re)generated handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-17 04:26:14"
$wsZhCn.Range("G3").Value = "2016-02-17 04:26:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-17 04:26:24"
$wsDeDe.Range("G3").Value = "2016-02-17 04:27:14"
